$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 33, shifting rows 33-35 down to 34-36
$ws.Rows.Item(33).Insert()

# Fill the new row 33 with data (copy of the pattern, new values)
$ws.Cells.Item(33, 1).Value = 5
$ws.Cells.Item(33, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(33, 3).Value = "Maule"
$ws.Cells.Item(33, 4).Value = 44491
$ws.Cells.Item(33, 4).NumberFormat = $ws.Cells.Item(34, 4).NumberFormat
$ws.Cells.Item(33, 5).Value = 7
$ws.Cells.Item(33, 6).Value = 300000000
$ws.Cells.Item(33, 7).Value = "Espárragos"
$ws.Cells.Item(33, 8).Value = "Verde"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 5000
$ws.Cells.Item(33, 11).Value = 850
$ws.Cells.Item(33, 12).Value = 850
$ws.Cells.Item(33, 13).Value = 850
$ws.Cells.Item(33, 14).Value = "`$/kilo"
$ws.Cells.Item(33, 15).Value = "Provincia de Linares"
$ws.Cells.Item(33, 16).Value = 850
$ws.Cells.Item(33, 17).Value = 1
$ws.Cells.Item(33, 18).Value = "Hortaliza"
